$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.25828618208765
$ws.Range("D2").Value = 7.523545576203079
$ws.Range("E2").Value = 12.99110612770389
$ws.Range("F2").Value = 39.77527920495806
$ws.Range("G2").Value = 46.55479462329361
$ws.Range("H2").Value = 18.75070558967343
$ws.Range("I2").Value = 29.97735987154836
$ws.Range("J2").Value = 10.43823966562165
$ws.Range("K2").Value = 12.14648907249011
$ws.Range("L2").Value = 11.11085322783178
$ws.Range("M2").Value = 15.89878456166634
$ws.Range("N2").Value = 22.17594656262283

$ws.Range("B3").Value = 14.17150529123666
$ws.Range("D3").Value = 7.525945128961435
$ws.Range("E3").Value = 12.99760025294545
$ws.Range("F3").Value = 39.83212891890038
$ws.Range("G3").Value = 46.57143555978646
$ws.Range("H3").Value = 18.79309673270426
$ws.Range("I3").Value = 30.06929591569883
$ws.Range("J3").Value = 10.4412275483072
$ws.Range("K3").Value = 11.95509246028224
$ws.Range("L3").Value = 11.10070362062641
$ws.Range("M3").Value = 15.88362711079523
$ws.Range("N3").Value = 22.23881252347295

$ws.Range("B4").Value = 14.12047337275954
$ws.Range("D4").Value = 7.527796075116091
$ws.Range("E4").Value = 13.0019315126586
$ws.Range("F4").Value = 39.87494313539423
$ws.Range("G4").Value = 46.59325810448418
$ws.Range("H4").Value = 18.82205844737025
$ws.Range("I4").Value = 30.12988288710011
$ws.Range("J4").Value = 10.44322981303387
$ws.Range("K4").Value = 11.83862762975129
$ws.Range("L4").Value = 11.09612419394336
$ws.Range("M4").Value = 15.87652329530162
$ws.Range("N4").Value = 22.27922534786825

$ws.Range("B5").Value = 14.10026149916168
$ws.Range("D5").Value = 7.528645602364032
$ws.Range("E5").Value = 13.00378313432594
$ws.Range("F5").Value = 39.89437684446278
$ws.Range("G5").Value = 46.60506661136633
$ws.Range("H5").Value = 18.83459782317892
$ws.Range("I5").Value = 30.15561350340877
$ws.Range("J5").Value = 10.44408800488055
$ws.Range("K5").Value = 11.79149390836675
$ws.Range("L5").Value = 11.09467553752227
$ws.Range("M5").Value = 15.8741854821578
$ws.Range("N5").Value = 22.29615107451317

$ws.Range("B6").Value = 14.09694111102983
$ws.Range("D6").Value = 7.528792427377727
$ws.Range("E6").Value = 13.00409582906589
$ws.Range("F6").Value = 39.89772371151072
$ws.Range("G6").Value = 46.60720339748706
$ws.Range("H6").Value = 18.83672449250954
$ws.Range("I6").Value = 30.15994893549382
$ws.Range("J6").Value = 10.44423306145323
$ws.Range("K6").Value = 11.7836889205778
$ws.Range("L6").Value = 11.09446025849184
$ws.Range("M6").Value = 15.87383102008241
$ws.Range("N6").Value = 22.2989892292204

$ws.Range("B7").Value = 14.12019840071297
$ws.Range("D7").Value = 7.527807146060899
$ws.Range("E7").Value = 13.00195613344656
$ws.Range("F7").Value = 39.87519718474719
$ws.Range("G7").Value = 46.59340555738447
$ws.Range("H7").Value = 18.82222457327041
$ws.Range("I7").Value = 30.13022568391865
$ws.Range("J7").Value = 10.4432412157119
$ws.Range("K7").Value = 11.83799056383989
$ws.Range("L7").Value = 11.09610296398403
$ws.Range("M7").Value = 15.87648950735998
$ws.Range("N7").Value = 22.27945176133209

$ws.Range("B8").Value = 14.22790767404657
$ws.Range("D8").Value = 7.52429474764928
$ws.Range("E8").Value = 12.99327401951866
$ws.Range("F8").Value = 39.79323855193365
$ws.Range("G8").Value = 46.55812277038749
$ws.Range("H8").Value = 18.76471301722367
$ws.Range("I8").Value = 30.0082008653273
$ws.Range("J8").Value = 10.43923513377651
$ws.Range("K8").Value = 12.080313906078
$ws.Range("L8").Value = 11.10701191261777
$ws.Range("M8").Value = 15.89310299474514
$ws.Range("N8").Value = 22.19724727799954

$ws.Range("B9").Value = 14.45617409665446
$ws.Range("D9").Value = 7.520390346497599
$ws.Range("E9").Value = 12.97897073429716
$ws.Range("F9").Value = 39.69533935501806
$ws.Range("G9").Value = 46.58105979796133
$ws.Range("H9").Value = 18.67522412872871
$ws.Range("I9").Value = 29.80172320712183
$ws.Range("J9").Value = 10.4327060984384
$ws.Range("K9").Value = 12.56116603464088
$ws.Range("L9").Value = 11.14142194139136
$ws.Range("M9").Value = 15.94301533269294
$ws.Range("N9").Value = 22.05036791108783

$ws.Range("B10").Value = 14.63312970831937
$ws.Range("D10").Value = 7.519323866086958
$ws.Range("E10").Value = 12.97011391172449
$ws.Range("F10").Value = 39.66178349627917
$ws.Range("G10").Value = 46.65404097967527
$ws.Range("H10").Value = 18.62369588474802
$ws.Range("I10").Value = 29.67000063016216
$ws.Range("J10").Value = 10.42871323265219
$ws.Range("K10").Value = 12.91432910929269
$ws.Range("L10").Value = 11.17450185630048
$ws.Range("M10").Value = 15.99004359113599
$ws.Range("N10").Value = 21.95110071316452

$ws.Range("B11").Value = 14.71538336379388
$ws.Range("D11").Value = 7.519226434080356
$ws.Range("E11").Value = 12.96644178576217
$ws.Range("F11").Value = 39.65485403633141
$ws.Range("G11").Value = 46.69938695428526
$ws.Range("H11").Value = 18.60334527773175
$ws.Range("I11").Value = 29.61440901344089
$ws.Range("J11").Value = 10.42707034281553
$ws.Range("K11").Value = 13.07423600092659
$ws.Range("L11").Value = 11.19121075089189
$ws.Range("M11").Value = 16.01363707003251
$ws.Range("N11").Value = 21.90780015590341

$ws.Range("B12").Value = 14.74676088670891
$ws.Range("D12").Value = 7.519244977977654
$ws.Range("E12").Value = 12.96510244298415
$ws.Range("F12").Value = 39.65342782781608
$ws.Range("G12").Value = 46.71829817258028
$ws.Range("H12").Value = 18.59608352197038
$ws.Range("I12").Value = 29.5939802129224
$ws.Range("J12").Value = 10.42647308894264
$ws.Range("K12").Value = 13.13462250082309
$ws.Range("L12").Value = 11.19777356914753
$ws.Range("M12").Value = 16.02288297208255
$ws.Range("N12").Value = 21.89166897245296

$ws.Range("B13").Value = 14.739993316494
$ws.Range("D13").Value = 7.519238523577229
$ws.Range("E13").Value = 12.96538861846867
$ws.Range("F13").Value = 39.65368173273158
$ws.Range("G13").Value = 46.71414806088333
$ws.Range("H13").Value = 18.59762769486267
$ws.Range("I13").Value = 29.59835223385091
$ws.Range("J13").Value = 10.42660061319575
$ws.Range("K13").Value = 13.12162552095791
$ws.Range("L13").Value = 11.19634972832206
$ws.Range("M13").Value = 16.02087792406744
$ws.Range("N13").Value = 21.89513130862704

$ws.Range("B14").Value = 14.71796033561315
$ws.Range("D14").Value = 7.519226850420576
$ws.Range("E14").Value = 12.96633057147972
$ws.Range("F14").Value = 39.6547127017474
$ws.Range("G14").Value = 46.70090797217662
$ws.Range("H14").Value = 18.60273893746696
$ws.Range("I14").Value = 29.61271584741797
$ws.Range("J14").Value = 10.42702070831198
$ws.Range("K14").Value = 13.07920770718502
$ws.Range("L14").Value = 11.19174597720453
$ws.Range("M14").Value = 16.0143915202152
$ws.Range("N14").Value = 21.906467713452

$ws.Range("B15").Value = 14.70449375518723
$ws.Range("D15").Value = 7.519226910737464
$ws.Range("E15").Value = 12.96691421098936
$ws.Range("F15").Value = 39.65550015871372
$ws.Range("G15").Value = 46.6930243555412
$ws.Range("H15").Value = 18.6059276248442
$ws.Range("I15").Value = 29.6215950468325
$ws.Range("J15").Value = 10.42728126569668
$ws.Range("K15").Value = 13.05320213861875
$ws.Range("L15").Value = 11.18895661779513
$ws.Range("M15").Value = 16.0104588382502
$ws.Range("N15").Value = 21.9134461687138

$ws.Range("B16").Value = 14.62778764141835
$ws.Range("D16").Value = 7.519338008585893
$ws.Range("E16").Value = 12.9703610659643
$ws.Range("F16").Value = 39.66240402804391
$ws.Range("G16").Value = 46.65132137792498
$ws.Range("H16").Value = 18.62508804427538
$ws.Range("I16").Value = 29.67372079112774
$ws.Range("J16").Value = 10.42882408397732
$ws.Range("K16").Value = 12.90385864700819
$ws.Range("L16").Value = 11.17344303209777
$ws.Range("M16").Value = 15.9885455683247
$ws.Range("N16").Value = 21.95396776482239

$ws.Range("B17").Value = 14.58116461196083
$ws.Range("D17").Value = 7.519505244610778
$ws.Range("E17").Value = 12.97256692939496
$ws.Range("F17").Value = 39.6687738412092
$ws.Range("G17").Value = 46.62884457892896
$ws.Range("H17").Value = 18.63763393239012
$ws.Range("I17").Value = 29.70680712537015
$ws.Range("J17").Value = 10.42981493369997
$ws.Range("K17").Value = 12.81200671510166
$ws.Range("L17").Value = 11.16434902074418
$ws.Range("M17").Value = 15.97566264612589
$ws.Range("N17").Value = 21.97930112716664

$ws.Range("B18").Value = 14.55451490724738
$ws.Range("D18").Value = 7.519637946972741
$ws.Range("E18").Value = 12.97386928197423
$ws.Range("F18").Value = 39.67322223037793
$ws.Range("G18").Value = 46.61706033471048
$ws.Range("H18").Value = 18.64514081518287
$ws.Range("I18").Value = 29.726245038787
$ws.Range("J18").Value = 10.43040117786005
$ws.Range("K18").Value = 12.75910883628157
$ws.Range("L18").Value = 11.159274906213
$ws.Range("M18").Value = 15.9684601171371
$ws.Range("N18").Value = 21.99404702770142

$ws.Range("B19").Value = 14.5455210772437
$ws.Range("D19").Value = 7.519689158278569
$ws.Range("E19").Value = 12.97431601003953
$ws.Range("F19").Value = 39.67486315461372
$ws.Range("G19").Value = 46.61326702594681
$ws.Range("H19").Value = 18.64773245627013
$ws.Range("I19").Value = 29.73289637788419
$ws.Range("J19").Value = 10.43060247763731
$ws.Range("K19").Value = 12.74118884063816
$ws.Range("L19").Value = 11.15758387727336
$ws.Range("M19").Value = 15.9660572256696
$ws.Range("N19").Value = 21.99906979793223

$ws.Range("B20").Value = 14.5861106303525
$ws.Range("D20").Value = 7.51948366540382
$ws.Range("E20").Value = 12.9723286349758
$ws.Range("F20").Value = 39.66801455987132
$ws.Range("G20").Value = 46.63111893648793
$ws.Range("H20").Value = 18.63626830030647
$ws.Range("I20").Value = 29.70324285134174
$ws.Range("J20").Value = 10.42970776613171
$ws.Range("K20").Value = 12.82179188761464
$ws.Range("L20").Value = 11.16530091761626
$ws.Range("M20").Value = 15.97701262538031
$ws.Range("N20").Value = 21.9765862627767

$ws.Range("B21").Value = 14.72442589992
$ws.Range("D21").Value = 7.519228777021909
$ws.Range("E21").Value = 12.96605250803083
$ws.Range("F21").Value = 39.65437738178337
$ws.Range("G21").Value = 46.70474975727949
$ws.Range("H21").Value = 18.60122557416582
$ws.Range("I21").Value = 29.60848001081868
$ws.Range("J21").Value = 10.42689664173747
$ws.Range("K21").Value = 13.09167183286675
$ws.Range("L21").Value = 11.19309184566209
$ws.Range("M21").Value = 16.01628831810688
$ws.Range("N21").Value = 21.90313073220548

$ws.Range("B22").Value = 14.81615197040081
$ws.Range("D22").Value = 7.519385187234845
$ws.Range("E22").Value = 12.9622491481731
$ws.Range("F22").Value = 39.65244574289806
$ws.Range("G22").Value = 46.76300743267556
$ws.Range("H22").Value = 18.58091453744396
$ws.Range("I22").Value = 29.55017562936882
$ws.Range("J22").Value = 10.42520435366046
$ws.Range("K22").Value = 13.26705619122902
$ws.Range("L22").Value = 11.21262596234529
$ws.Range("M22").Value = 16.04377087332548
$ws.Range("N22").Value = 21.85667203162186

$ws.Range("B23").Value = 14.76708197895713
$ws.Range("D23").Value = 7.519272258683912
$ws.Range("E23").Value = 12.96425179971213
$ws.Range("F23").Value = 39.65283834949908
$ws.Range("G23").Value = 46.73098957283273
$ws.Range("H23").Value = 18.59151773713232
$ws.Range("I23").Value = 29.58096177225329
$ws.Range("J23").Value = 10.42609432150565
$ws.Range("K23").Value = 13.17356046978804
$ws.Range("L23").Value = 11.20207591438533
$ws.Range("M23").Value = 16.02893862485119
$ws.Range("N23").Value = 21.88132659119717

$ws.Range("B24").Value = 14.58387405247928
$ws.Range("D24").Value = 7.519493307493977
$ws.Range("E24").Value = 12.97243626145532
$ws.Range("F24").Value = 39.66835538160755
$ws.Range("G24").Value = 46.63008715418754
$ws.Range("H24").Value = 18.6368847865085
$ws.Range("I24").Value = 29.70485296335401
$ws.Range("J24").Value = 10.42975616491087
$ws.Range("K24").Value = 12.81736828995986
$ws.Range("L24").Value = 11.16487008461279
$ws.Range("M24").Value = 15.976401663844
$ws.Range("N24").Value = 21.97781308776393

$ws.Range("B25").Value = 14.39271244872119
$ws.Range("D25").Value = 7.521128906549457
$ws.Range("E25").Value = 12.9825495234807
$ws.Range("F25").Value = 39.71508838165555
$ws.Range("G25").Value = 46.56498846628491
$ws.Range("H25").Value = 18.69693744314867
$ws.Range("I25").Value = 29.85407091560747
$ws.Range("J25").Value = 10.43433083176724
$ws.Range("K25").Value = 12.43085002537498
$ws.Range("L25").Value = 11.13073278361133
$ws.Range("M25").Value = 15.92767758453951
$ws.Range("N25").Value = 22.08857827405298

